$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") updates ---
# These source values are plain numeric-looking text (e.g. "231.54") that must
# remain *text* cells (matching the original inlineStr storage), not be auto-
# converted to numbers by Excel. Prefixing with a literal apostrophe forces
# text entry; ClearFormats() then strips the transient quote-prefix style Excel
# applies so the cell keeps its original (unstyled) appearance.
$ws.Range("D2").Value = "'28.658.48"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").Value = "'1.798.48"
$ws.Range("D3").ClearFormats()
$ws.Range("D5").Value = "'231.54"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").Value = "'0.5888"
$ws.Range("D6").ClearFormats()
$ws.Range("D8").Value = "'0.2767"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").Value = "'0.06802"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").Value = "'23.21"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").Value = "'0.07533"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'1.819.86"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "'4.781"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").Value = "'0.6197"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Value = "'2.043.38"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = "'0.000009118"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").Value = "'75.64"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").Value = "'28.645.28"
$ws.Range("D18").ClearFormats()
$ws.Range("D19").Value = "'5.473"
$ws.Range("D19").ClearFormats()
$ws.Range("D21").Value = "'210.81"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").Value = "'11.51"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").Value = "'6.823"
$ws.Range("D23").ClearFormats()
$ws.Range("D24").Value = "'1.005"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").Value = "'153.63"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").Value = "'7.907"
$ws.Range("D26").ClearFormats()
$ws.Range("D27").Value = "'0.1268"
$ws.Range("D27").ClearFormats()
$ws.Range("D28").Value = "'16.44"
$ws.Range("D28").ClearFormats()
$ws.Range("D29").Value = "'1.428"
$ws.Range("D29").ClearFormats()
$ws.Range("D30").Value = "'0.06117"
$ws.Range("D30").ClearFormats()
$ws.Range("D32").Value = "'3.812"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").Value = "'3.784"
$ws.Range("D33").ClearFormats()
$ws.Range("D34").Value = "'1.735"
$ws.Range("D34").ClearFormats()
$ws.Range("D35").Value = "'1.055"
$ws.Range("D35").ClearFormats()
$ws.Range("D36").Value = "'0.6428"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").Value = "'2.497"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").Value = "'2.716"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").Value = "'6.542"
$ws.Range("D39").ClearFormats()
$ws.Range("D40").Value = "'0.01697"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").Value = "'1.147.93"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").Value = "'0.8869"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").Value = "'1.006"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").Value = "'100.20"
$ws.Range("D44").ClearFormats()
$ws.Range("D46").Value = "'60.21"
$ws.Range("D46").ClearFormats()
$ws.Range("D47").Value = "'0.00000000112"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").Value = "'1.591"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").Value = "'8.351"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").Value = "'0.05470"
$ws.Range("D50").ClearFormats()

# --- Column E ("Volume(1h)") updates ---
# These are padded percentage strings (e.g. "  -2.03%  ") which Excel never
# interprets as numbers because of the surrounding whitespace, so a plain
# text assignment is safe and keeps the default cell style untouched.
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("E9").Value = "  -3.24%  "
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("E16").Value = "  -7.98%  "
$ws.Range("E17").Value = "  -4.25%  "
$ws.Range("E18").Value = "  -1.97%  "
$ws.Range("E19").Value = "  -6.05%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").Value = "  -5.79%  "
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("E23").Value = "  -2.46%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("E29").Value = "  -3.59%  "
$ws.Range("E30").Value = "  -1.35%  "
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("E35").Value = "  -5.41%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -1.84%  "
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("E41").Value = "  -5.95%  "
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("E45").Value = "  -1.84%  "
$ws.Range("E46").Value = "  -3.45%  "
$ws.Range("E47").Value = "  -3.52%  "
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("E49").Value = "  -1.97%  "
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("E51").Value = "  -1.77%  "
